# RCP-Qlearning now matches the theoretical result
# 1. fix an error in RTT-Variance calculation
# 2. fix corresponding errors in the theoretical calculation because of the error in 1.
# 3. fix the mistake of failed to normalize the delay when calculating expected delay
#
# The stale RCP-Q_Learning simulation output (range "Q_learning_v_1",
# case_study!$AB$3:$AK$10) is cleared out for the rows that need to be
# regenerated with the corrected calculation, mirroring the author's
# re-run/re-import of Q_learning_v.txt after fixing the upstream bug.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("case_study")

# --- Clear the now-invalid simulation numbers (rows 5-10, cols AB:AK) ---
# Row 5: keep the last two columns (AJ/AK), clear the rest.
$ws.Range("AB5:AI5").ClearContents()
# Rows 6-9: keep only the final column (AK), clear the rest.
$ws.Range("AB6:AJ6").ClearContents()
$ws.Range("AB7:AJ7").ClearContents()
$ws.Range("AB8:AJ8").ClearContents()
$ws.Range("AB9:AJ9").ClearContents()
# Row 10: clear the entire row of data (nothing survives recomputation yet).
$ws.Range("AB10:AK10").ClearContents()

# --- Move the sheet's active selection, as left by the author after editing ---
$ws.Range("AD29").Select() | Out-Null
